# VOLAIB.xlsx update: refresh EffectiveDate / PreviousExpDate values
# from 07302023 to 08302023 on the aibcustomerInfo sheet ("MTR366 added").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "08302023"
$ws.Range("I2").Value = "08302023"
$ws.Range("F3").Value = "08302023"
$ws.Range("I3").Value = "08302023"

# Leave the selection where the edit ended, matching the saved cursor position.
$ws.Range("E6").Select()
